# "opt for stoplose and holiday"
# - add a new "beta" column (U) with per-symbol values
# - move the old stray note text out of U2 into new footnote rows (27, 29, 30)
# - change the 4H row's (row 24) timeframe markers from "1h" to "4h"
# - unhide columns N:Q
# - update the active selection shown in the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- unhide columns N:Q (14-17) ---
$ws.Columns.Item(14).Hidden = $False
$ws.Columns.Item(15).Hidden = $False
$ws.Columns.Item(16).Hidden = $False
$ws.Columns.Item(17).Hidden = $False

# --- new "beta" header + per-row values ---
$ws.Range("U1").Value = "beta"
$ws.Range("U2").Value = 622
$ws.Range("U4").Value = 655
$ws.Range("U7").Value = 442
$ws.Range("U9").Value = -8
$ws.Range("U13").Value = 603
$ws.Range("U22").Value = 1779
$ws.Range("U24").Value = 500

# --- row 24 (4H) stoploss timeframe fix: 1h -> 4h ---
$ws.Range("P24").Value = "4h"
$ws.Range("R24").Value = "4h"

# --- footnotes moved below the table ---
$ws.Range("A27").Value = "注:"
$ws.Range("A29").Value = "USD,CAD经济体关系大，所以EURUSD与EURCAD共用仓位额度"
$ws.Range("A30").Value = "非农数据公布天强制close会造成利润大幅减少，建议4H级别以上不强制close，只是不开新单"

# --- update the view's stored selection ---
[void]$ws.Range("B1:M1048576").Select()

Write-Output "done"
